$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '255.18'
$c.ClearFormats()

$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '3.98%'
$c.ClearFormats()

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '27.65'
$c.ClearFormats()

$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '-7.50%'
$c.ClearFormats()

$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '5.186'
$c.ClearFormats()

$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '0.43%'
$c.ClearFormats()

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '0.05859'
$c.ClearFormats()

$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '1.97%'
$c.ClearFormats()

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '6.718'
$c.ClearFormats()

$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '1.11%'
$c.ClearFormats()

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.8685'
$c.ClearFormats()

$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '1.15%'
$c.ClearFormats()

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.9489'
$c.ClearFormats()

$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '11.29%'
$c.ClearFormats()

$c = $ws.Range('B9')
$c.NumberFormat = '@'
$c.Value = 'One'
$c.ClearFormats()

$c = $ws.Range('C9')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$c.ClearFormats()

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.0006078'
$c.ClearFormats()

$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '2.36%'
$c.ClearFormats()

$c = $ws.Range('B10')
$c.NumberFormat = '@'
$c.Value = 'WazirX'
$c.ClearFormats()

$c = $ws.Range('C10')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$c.ClearFormats()

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.1408'
$c.ClearFormats()

$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '1.92%'
$c.ClearFormats()

$c = $ws.Range('B11')
$c.NumberFormat = '@'
$c.Value = 'MandalaExchangeToken'
$c.ClearFormats()

$c = $ws.Range('C11')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$c.ClearFormats()

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.07165'
$c.ClearFormats()

$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '1.15%'
$c.ClearFormats()

$c = $ws.Range('B12')
$c.NumberFormat = '@'
$c.Value = 'BitrueCoin'
$c.ClearFormats()

$c = $ws.Range('C12')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$c.ClearFormats()

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.03179'
$c.ClearFormats()

$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '1.20%'
$c.ClearFormats()

$c = $ws.Range('B13')
$c.NumberFormat = '@'
$c.Value = 'BitMartToken'
$c.ClearFormats()

$c = $ws.Range('C13')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$c.ClearFormats()

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.09232'
$c.ClearFormats()

$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '-1.49%'
$c.ClearFormats()

$c = $ws.Range('B14')
$c.NumberFormat = '@'
$c.Value = 'BitForexToken'
$c.ClearFormats()

$c = $ws.Range('C14')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$c.ClearFormats()

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.001548'
$c.ClearFormats()

$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '1.56%'
$c.ClearFormats()

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.006001'
$c.ClearFormats()

$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '-1.66%'
$c.ClearFormats()

$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '-0.87%'
$c.ClearFormats()

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '3.226'
$c.ClearFormats()

$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '-1.79%'
$c.ClearFormats()

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '2.204'
$c.ClearFormats()

$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '0.59%'
$c.ClearFormats()

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.3174'
$c.ClearFormats()

$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '0.64%'
$c.ClearFormats()

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.03442'
$c.ClearFormats()

$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '3.86%'
$c.ClearFormats()

$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '-0.34%'
$c.ClearFormats()

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '3.532'
$c.ClearFormats()

$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '1.29%'
$c.ClearFormats()

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.04177'
$c.ClearFormats()

$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '1.48%'
$c.ClearFormats()

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.001226'
$c.ClearFormats()

$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '0.48%'
$c.ClearFormats()

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.004791'
$c.ClearFormats()

$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '15.16%'
$c.ClearFormats()

$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '-0.03%'
$c.ClearFormats()

$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '1.19%'
$c.ClearFormats()

$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '1.56%'
$c.ClearFormats()

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.005668'
$c.ClearFormats()

$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '58.27%'
$c.ClearFormats()

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.1102'
$c.ClearFormats()

$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '2.79%'
$c.ClearFormats()

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.002299'
$c.ClearFormats()

$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '-6.53%'
$c.ClearFormats()

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.009804'
$c.ClearFormats()

$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '-1.63%'
$c.ClearFormats()

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.00005366'
$c.ClearFormats()

$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '-1.74%'
$c.ClearFormats()

$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '0.00%'
$c.ClearFormats()

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.09997'
$c.ClearFormats()

$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '11.25%'
$c.ClearFormats()

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.002128'
$c.ClearFormats()

$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '-4.00%'
$c.ClearFormats()

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.00002099'
$c.ClearFormats()

$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '0.00%'
$c.ClearFormats()

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.0001999'
$c.ClearFormats()

$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '0.00%'
$c.ClearFormats()
